# Fix error in current CLO size calculation: the underlying distribution
# summary statistics now only add up offered tranches, which shifts the
# histogram bin edges (the bracketed "[lo-hi]" labels in column A) on each
# of the three frequency sheets, and moves a couple of the "1" frequency
# hits to the adjacent bin on the "Proj Equity Yield" and "WA Adv Rate"
# sheets.

$wb = $excel.ActiveWorkbook

# --- "WA Cost of Funds" sheet: refreshed bin-edge labels (frequency counts unchanged) ---
$wsCOF = $wb.Worksheets.Item("WA Cost of Funds")
$cofLabels = @(
    "[3.539-3.547]",
    "[3.547-3.556]",
    "[3.556-3.564]",
    "[3.564-3.573]",
    "[3.573-3.581]",
    "[3.581-3.59]",
    "[3.59-3.598]",
    "[3.598-3.607]",
    "[3.607-3.615]",
    "[3.615-3.623]",
    "[3.623-3.632]",
    "[3.632-3.64]",
    "[3.64-3.649]",
    "[3.649-3.657]"
)
for ($i = 0; $i -lt $cofLabels.Length; $i++) {
    $wsCOF.Range("A" + ($i + 1)).Value = $cofLabels[$i]
}

# --- "Proj Equity Yield" sheet: refreshed bin-edge labels + one frequency hit moved ---
$wsEY = $wb.Worksheets.Item("Proj Equity Yield")
$eyLabels = @(
    "[12.3-12.34]",
    "[12.34-12.38]",
    "[12.38-12.42]",
    "[12.42-12.46]",
    "[12.46-12.5]",
    "[12.5-12.54]",
    "[12.54-12.58]",
    "[12.58-12.62]",
    "[12.62-12.66]",
    "[12.66-12.7]",
    "[12.7-12.74]",
    "[12.74-12.78]",
    "[12.78-12.82]",
    "[12.82-12.86]"
)
for ($i = 0; $i -lt $eyLabels.Length; $i++) {
    $wsEY.Range("A" + ($i + 1)).Value = $eyLabels[$i]
}
$wsEY.Range("B12").Value = 0
$wsEY.Range("B13").Value = 1

# --- "WA Adv Rate" sheet: refreshed bin-edge labels + two frequency hits moved ---
$wsAR = $wb.Worksheets.Item("WA Adv Rate")
$arLabels = @(
    "[81.83-81.84]",
    "[81.84-81.85]",
    "[81.85-81.86]",
    "[81.86-81.86]",
    "[81.86-81.87]",
    "[81.87-81.88]",
    "[81.88-81.89]",
    "[81.89-81.89]",
    "[81.89-81.9]",
    "[81.9-81.91]",
    "[81.91-81.92]",
    "[81.92-81.92]",
    "[81.92-81.93]",
    "[81.93-81.94]"
)
for ($i = 0; $i -lt $arLabels.Length; $i++) {
    $wsAR.Range("A" + ($i + 1)).Value = $arLabels[$i]
}
$wsAR.Range("B1").Value = 0
$wsAR.Range("B2").Value = 1
$wsAR.Range("B13").Value = 1
$wsAR.Range("B14").Value = 0
